# Commit: "enhancement - The TreeView is now updated as soon as a new link
# is created." — underlying OOXML change: the template's fixed
# Header&Footer date stamp moves from 11/07/15 to 12/07/15 (applied on the
# slide master and every slide layout), and three connected shapes in the
# class-diagram slide ("CodeEditor" rounded rectangle plus its two bent
# connectors) are nudged/resized so the diagram lines reconnect cleanly.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fixed date placeholder: 11/07/15 -> 12/07/15 on the slide master and
#    on every slide layout (mirrors PowerPoint's "Apply to All" on the
#    Insert > Header & Footer dialog with a fixed date).
# ---------------------------------------------------------------------
function Update-FixedDate($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "11/07/15") {
                    $tr.Text = "12/07/15"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-FixedDate $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-FixedDate $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Reposition/resize the "CodeEditor" shape and the two connectors that
#    link it into the diagram, on slide 1.
# ---------------------------------------------------------------------
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$EMUPerPoint = 12700.0
$slide = $p.Slides.Item(1)

# "Afgeronde rechthoek 19" (CodeEditor rounded rectangle) - move left edge.
$codeEditor = Get-ShapeById $slide.Shapes 49
$codeEditor.Left = 7292806 / $EMUPerPoint

# "Elbow Connector 65" feeding into the CodeEditor shape - shorten it.
$elbow = Get-ShapeById $slide.Shapes 66
$elbow.Width = 358607 / $EMUPerPoint

# "Shape 68" bent connector above the CodeEditor shape - shorten it.
$bentConn = Get-ShapeById $slide.Shapes 69
$bentConn.Width = 1071768 / $EMUPerPoint
